# Weekly update: insert 3 new price rows for "Agrícola del Norte S.A. de Arica - Locoto"
# 1) a new row right before the current row 55 (pushes the old 55..90 block down by one)
# 2) two brand-new rows appended right after the (now shifted) old row 89, i.e. at 91/92,
#    which pushes the old row 90 (now at 91) further down to 93.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: insert a single row at position 55 ---------------------------
$ws.Rows.Item(55).Insert()

$ws.Cells.Item(55, 1).Value = 1
$ws.Cells.Item(55, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(55, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(55, 4).Value = 44567
$ws.Cells.Item(55, 5).Value = 15
$ws.Cells.Item(55, 6).Value = 100112042
$ws.Cells.Item(55, 7).Value = "Locoto"
$ws.Cells.Item(55, 8).Value = "Sin especificar"
$ws.Cells.Item(55, 9).Value = "Primera"
$ws.Cells.Item(55, 10).Value = 120
$ws.Cells.Item(55, 11).Value = 14000
$ws.Cells.Item(55, 12).Value = 15000
$ws.Cells.Item(55, 13).Value = 14500
$ws.Cells.Item(55, 14).Value = "$/caja 20 kilos"
$ws.Cells.Item(55, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(55, 16).Value = 725
$ws.Cells.Item(55, 17).Value = 20
$ws.Cells.Item(55, 18).Value = "Hortaliza"

# --- Step 2: insert two new rows at position 91 ----------------------------
$ws.Rows.Item(91).Insert()
$ws.Rows.Item(91).Insert()

$ws.Cells.Item(91, 1).Value = 1
$ws.Cells.Item(91, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(91, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(91, 4).Value = 44568
$ws.Cells.Item(91, 5).Value = 15
$ws.Cells.Item(91, 6).Value = 100112042
$ws.Cells.Item(91, 7).Value = "Locoto"
$ws.Cells.Item(91, 8).Value = "Sin especificar"
$ws.Cells.Item(91, 9).Value = "Primera"
$ws.Cells.Item(91, 10).Value = 120
$ws.Cells.Item(91, 11).Value = 14000
$ws.Cells.Item(91, 12).Value = 15000
$ws.Cells.Item(91, 13).Value = 14500
$ws.Cells.Item(91, 14).Value = "$/caja 20 kilos"
$ws.Cells.Item(91, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(91, 16).Value = 725
$ws.Cells.Item(91, 17).Value = 20
$ws.Cells.Item(91, 18).Value = "Hortaliza"

$ws.Cells.Item(92, 1).Value = 1
$ws.Cells.Item(92, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(92, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(92, 4).Value = 44568
$ws.Cells.Item(92, 5).Value = 15
$ws.Cells.Item(92, 6).Value = 100112042
$ws.Cells.Item(92, 7).Value = "Locoto"
$ws.Cells.Item(92, 8).Value = "Sin especificar"
$ws.Cells.Item(92, 9).Value = "Segunda"
$ws.Cells.Item(92, 10).Value = 120
$ws.Cells.Item(92, 11).Value = 10000
$ws.Cells.Item(92, 12).Value = 11000
$ws.Cells.Item(92, 13).Value = 10500
$ws.Cells.Item(92, 14).Value = "$/caja 20 kilos"
$ws.Cells.Item(92, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(92, 16).Value = 525
$ws.Cells.Item(92, 17).Value = 20
$ws.Cells.Item(92, 18).Value = "Hortaliza"
